$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record for row 236 needs to be inserted; the existing row 236
# (and everything below it through row 278) shifts down by one row.
$ws.Rows("237:237").Insert()

# The row that used to be at 236 now needs to be duplicated into the freshly
# inserted row 237 (its contents are unchanged from before the edit).
$ws.Range("A236:R236").Copy()
$ws.Range("A237").PasteSpecial()

# Row 236 itself becomes the new weekly entry: same market/category/quality
# metadata, but a new date and new price figures.
$ws.Range("D236").Value = 44722
$ws.Range("K236").Value = 7000
$ws.Range("L236").Value = 7500
$ws.Range("M236").Value = 7227
$ws.Range("P236").Value = 1204
